# Resolved errors and added testng xml files for every test classes
#
# The underlying data fix: on the "Negative Testdata" sheet a couple of
# stray/incorrect test values are replaced with the correct ones.

$wb = $excel.ActiveWorkbook

$wsNeg = $wb.Worksheets.Item("Negative Testdata")

# B3 held the stray password "df@456" - fix it to the valid "Test@123"
$wsNeg.Range("B3").Value = "Test@123"

# A4 held the stray email "ghi@g.com" - fix it to the valid "sk@g.com"
$wsNeg.Range("A4").Value = "sk@g.com"

# Make "Negative Testdata" the active sheet and leave the selection on C9,
# matching the cursor position the author had when the file was saved.
$wsNeg.Activate()
$wsNeg.Range("C9").Select()

$wb.Save()
